# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
#
# Updates the td_sim_1 (column C) and record_atd (column D) values for
# rows 2-18, and the average_simulation_TD value in C19, on the active
# worksheet (Sheet1) of tds_sim.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (new C value, new D value). D is omitted (left unchanged) for
# rows where the diff does not touch column D (none here besides row 19,
# which has no D cell at all).
$updates = @{
    2  = @{ C = 20;    D = 16 }
    3  = @{ C = 15;    D = 13 }
    4  = @{ C = 38;    D = 35 }
    5  = @{ C = 5;     D = 4 }
    6  = @{ C = 44;    D = 41.5 }
    7  = @{ C = 35;    D = 35 }
    8  = @{ C = 62;    D = 58 }
    9  = @{ C = 48;    D = 44 }
    10 = @{ C = 18;    D = 13 }
    11 = @{ C = 6;     D = 3.5 }
    12 = @{ C = 29;    D = 26 }
    13 = @{ C = 60;    D = 64 }
    14 = @{ C = 9;     D = 9 }
    15 = @{ C = 16;    D = 13 }
    16 = @{ C = 7;     D = 5.5 }
    17 = @{ C = 148;   D = 152.5 }
    18 = @{ C = 21;    D = 17 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
}

# Row 19 only has a new value for C (the overall average); there is no D19.
$ws.Range("C19").Value = 34.1764705882353
